$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) that switch from the
#    deck's custom "Table_0" style to the built-in table style
#    {05B0D78D-179F-48DA-B335-512BDB2C7FE8}.
# ---------------------------------------------------------------------------
$newTableStyleId = "{05B0D78D-179F-48DA-B335-512BDB2C7FE8}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme color palette back to the stock
#    "Office" colors (the deck currently carries the custom "Red Violet"
#    / Integral palette on the live theme).
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501     # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407       # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456     # accent6  70AD47
    11 = 12673797    # hlink    0563C1
    12 = 7491477     # folHlink 954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i]
}
